# feat: add 2022-Q1 data
#
# The workbook has quarterly "holdings detail" sheets (2020-Q4 .. 2021-Q4)
# plus a running "总计" (totals) summary sheet. This change:
#   1. Adds a new "2022-Q1" detail sheet (same shape as the other quarter
#      sheets) right before "总计". It takes over 总计's old sheetId/rId
#      slot -- matching how the authors actually produced this commit: the
#      old "总计" sheet object was renamed/repurposed into "2022-Q1" and a
#      fresh "总计" sheet was appended right after it.
#   2. The (new) "总计" sheet gets a new first data row for 2022-Q1
#      (4 holdings, 0.30 亿元), with the older rows pushed down.

$wb = $excel.ActiveWorkbook

$quarterSheet  = $wb.Worksheets.Item("2021-Q4")  # template to copy formatting from
$oldTotalSheet = $wb.Worksheets.Item("总计")

function Set-TextCell($rng, $text) {
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# --- Step 1: repurpose the *old* 总计 sheet object into the new "2022-Q1"
# detail sheet, then append the new blank "总计" sheet right after it.
$detailSheet = $oldTotalSheet
$detailSheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Add($null, $detailSheet)
$totalSheet.Name = "总计"

# =====================================================================
# Step 2: build the "2022-Q1" detail sheet (fund holdings, same shape as
# the other quarterly detail sheets)
# =====================================================================

# clear out the old 总计 rows/columns that used to live here
$detailSheet.Range("A1:D10").Clear()

# header row, with formatting copied from the 2021-Q4 template sheet
$quarterSheet.Range("B1:H1").Copy()
$detailSheet.Range("B1").PasteSpecial(-4122)

$detailSheet.Range("B1").Value = "基金代码"
$detailSheet.Range("C1").Value = "基金名称"
$detailSheet.Range("D1").Value = "基金规模"
$detailSheet.Range("E1").Value = "股票总仓位"
$detailSheet.Range("F1").Value = "仓位占比"
$detailSheet.Range("G1").Value = "持有市值(亿元)"
$detailSheet.Range("H1").Value = "仓位排名"

# index column (A2:A5) formatting, copied from the template sheet
$quarterSheet.Range("A2:A5").Copy()
$detailSheet.Range("A2").PasteSpecial(-4122)

# row 2: 001304 建信鑫安回报灵活配置混合
$detailSheet.Range("A2").Value = 0
Set-TextCell $detailSheet.Range("B2") "001304"
Set-TextCell $detailSheet.Range("C2") "建信鑫安回报灵活配置混合"
Set-TextCell $detailSheet.Range("D2") "2.13"
Set-TextCell $detailSheet.Range("E2") "66.83"
Set-TextCell $detailSheet.Range("F2") "5.53"
Set-TextCell $detailSheet.Range("G2") "0.1178"
$detailSheet.Range("H2").Value = 8

# row 3: 006279 中金瑞祥灵活配置混合A
$detailSheet.Range("A3").Value = 1
Set-TextCell $detailSheet.Range("B3") "006279"
Set-TextCell $detailSheet.Range("C3") "中金瑞祥灵活配置混合A"
Set-TextCell $detailSheet.Range("D3") "2.10"
Set-TextCell $detailSheet.Range("E3") "59.54"
Set-TextCell $detailSheet.Range("F3") "4.78"
Set-TextCell $detailSheet.Range("G3") "0.1004"
$detailSheet.Range("H3").Value = 7

# row 4: 005396 中金丰硕混合
$detailSheet.Range("A4").Value = 2
Set-TextCell $detailSheet.Range("B4") "005396"
Set-TextCell $detailSheet.Range("C4") "中金丰硕混合"
Set-TextCell $detailSheet.Range("D4") "1.83"
Set-TextCell $detailSheet.Range("E4") "71.47"
Set-TextCell $detailSheet.Range("F4") "4.70"
Set-TextCell $detailSheet.Range("G4") "0.0860"
$detailSheet.Range("H4").Value = 10

# row 5: 006280 中金瑞祥灵活配置混合C
$detailSheet.Range("A5").Value = 3
Set-TextCell $detailSheet.Range("B5") "006280"
Set-TextCell $detailSheet.Range("C5") "中金瑞祥灵活配置混合C"
Set-TextCell $detailSheet.Range("D5") "0.00"
Set-TextCell $detailSheet.Range("E5") "59.54"
Set-TextCell $detailSheet.Range("F5") "4.78"
$detailSheet.Range("G5").NumberFormat = "General"
$detailSheet.Range("G5").Value = 0
$detailSheet.Range("H5").Value = 7

# =====================================================================
# Step 3: build the (new, blank) "总计" sheet: header + the 2022-Q1 row
# on top, followed by the pre-existing quarters (unchanged values),
# each row's index column renumbered to match its new position.
# =====================================================================

# header row, with formatting copied from the detail sheet's header style
$detailSheet.Range("B1:D1").Copy()
$totalSheet.Range("B1").PasteSpecial(-4122)
$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

# index column (A2:A7) formatting, copied from the template sheet
$quarterSheet.Range("A2:A7").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalData = @(
    @(0, "2022-Q1", 4, 0.3),
    @(1, "2021-Q4", 9, 0.76),
    @(2, "2021-Q3", 3, 0.23),
    @(3, "2021-Q2", 2, 0.08),
    @(4, "2021-Q1", 3, 0.06),
    @(5, "2020-Q4", 1, 0.24)
)

$r = 2
foreach ($row in $totalData) {
    $totalSheet.Range("A$r").Value = $row[0]
    $totalSheet.Range("B$r").Value = $row[1]
    $totalSheet.Range("C$r").Value = $row[2]
    $totalSheet.Range("D$r").Value = $row[3]
    $r = $r + 1
}
